{"js": "// Office.js (Word JavaScript API) reproduction of the commit:\n// \"Minor update and changes (Updated todo and fixed typo, changed keyword)\"\n//\n// Summary of the change being applied (see diff):\n//  1. Paragraph \"Umwandlung zu PNAS (Andre), wenn m\u00f6glich mit Harvard\n//     Citation Stil\" -> \"Umwandlung zu PNAS,wenn m\u00f6glich mit Harvard\n//     Citation Stil (Andre),\" (\"(Andre)\" moved to the end) with\n//     \"Umwandlung zu PNAS\" colored green (#00B050) and\n//     \"wenn m\u00f6glich mit Harvard Citation Stil \" colored amber (#FFC000).\n//  2. Paragraph \"Interpretation (Andre)\" loses the `_GoBack` bookmark\n//     that used to sit inside it (text itself is unchanged).\n//  3. Paragraph starting \"Muss: ...\" gets two of its spans colored green\n//     (#00B050): \"Titel, ... (warum ist es wichtig), \" and\n//     \"Literaturangabe\" (text itself is unchanged).\n//  4. The `_GoBack` bookmark re-appears on the empty paragraph right\n//     after the \"Kann: ...\" paragraph (near \"1. Einleitung\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) \"Umwandlung zu PNAS ...\" paragraph: reorder \"(Andre)\" to the end\n//    and recolor the pieces.\n// ---------------------------------------------------------------------\nlet umwandlungPara = null;\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Umwandlung zu PNAS\") !== -1) {\n    umwandlungPara = p;\n    break;\n  }\n}\n\nif (umwandlungPara) {\n  umwandlungPara.insertText(\n    \"Umwandlung zu PNAS,wenn m\u00f6glich mit Harvard Citation Stil (Andre),\",\n    \"Replace\"\n  );\n  await context.sync();\n\n  const greenRange = umwandlungPara\n    .getRange()\n    .search(\"Umwandlung zu PNAS\", { matchCase: true });\n  greenRange.load(\"items\");\n  await context.sync();\n  if (greenRange.items.length > 0) {\n    greenRange.items[0].font.color = \"#00B050\";\n  }\n\n  const amberRange = umwandlungPara\n    .getRange()\n    .search(\"wenn m\u00f6glich mit Harvard Citation Stil \", { matchCase: true });\n  amberRange.load(\"items\");\n  await context.sync();\n  if (amberRange.items.length > 0) {\n    amberRange.items[0].font.color = \"#FFC000\";\n  }\n\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) \"Interpretation (Andre)\" paragraph: remove the stray `_GoBack`\n//    bookmark (its text does not change) -- deleteBookmark removes only\n//    the bookmark markers, not the underlying text.\n// ---------------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) \"Muss: ...\" paragraph: color two spans green, text unchanged.\n// ---------------------------------------------------------------------\nlet mussPara = null;\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\nfor (const p of paragraphs2.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs2.items) {\n  if (p.text.indexOf(\"Muss:\") === 0) {\n    mussPara = p;\n    break;\n  }\n}\n\nif (mussPara) {\n  const greenSpan1 = mussPara\n    .getRange()\n    .search(\n      \"Titel, Autorenliste, Affiliationen, Abstract (Englisch), Einleitung (\\\\dropcap) (warum ist es wichtig), \",\n      { matchCase: true }\n    );\n  greenSpan1.load(\"items\");\n  await context.sync();\n  if (greenSpan1.items.length > 0) {\n    greenSpan1.items[0].font.color = \"#00B050\";\n  }\n\n  const greenSpan2 = mussPara\n    .getRange()\n    .search(\"Literaturangabe\", { matchCase: true });\n  greenSpan2.load(\"items\");\n  await context.sync();\n  if (greenSpan2.items.length > 0) {\n    greenSpan2.items[0].font.color = \"#00B050\";\n  }\n\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 4) Re-insert the `_GoBack` bookmark on the empty paragraph right\n//    after \"Kann: ...\" (Word moves its \"last edit\" bookmark here).\n// ---------------------------------------------------------------------\nlet kannPara = null;\nconst paragraphs3 = body.paragraphs;\nparagraphs3.load(\"items\");\nawait context.sync();\nfor (const p of paragraphs3.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs3.items.length; i++) {\n  if (paragraphs3.items[i].text.indexOf(\"Kann:\") === 0) {\n    kannPara = paragraphs3.items[i + 1];\n    break;\n  }\n}\n\nif (kannPara) {\n  kannPara.getRange().insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop reproduction of the commit:\n# \"Minor update and changes (Updated todo and fixed typo, changed keyword)\"\n#\n# Summary of the change being applied (see diff):\n#  1. Paragraph \"Umwandlung zu PNAS (Andre), wenn m\u00f6glich mit Harvard\n#     Citation Stil\" -> \"Umwandlung zu PNAS,wenn m\u00f6glich mit Harvard\n#     Citation Stil (Andre),\" (\"(Andre)\" moved to the end) with\n#     \"Umwandlung zu PNAS\" colored green (#00B050) and\n#     \"wenn m\u00f6glich mit Harvard Citation Stil \" colored amber (#FFC000).\n#  2. Paragraph \"Interpretation (Andre)\" loses the `_GoBack` bookmark\n#     that used to sit inside it (text itself is unchanged).\n#  3. Paragraph starting \"Muss: ...\" gets two of its spans colored green\n#     (#00B050): \"Titel, ... (warum ist es wichtig), \" and\n#     \"Literaturangabe\" (text itself is unchanged).\n#  4. The `_GoBack` bookmark re-appears on the empty paragraph right\n#     after the \"Kann: ...\" paragraph (near \"1. Einleitung\").\n\n$d = $word.ActiveDocument\n\n# Word's Font.Color is a BGR-packed integer (0x00BBGGRR), not the RRGGBB\n# hex used in OOXML <w:color w:val=\"RRGGBB\"/>.\nfunction RgbToWordColor([string]$hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    return ($b * 65536) + ($g * 256) + $r\n}\n\n$green = RgbToWordColor(\"00B050\")\n$amber = RgbToWordColor(\"FFC000\")\n\n# ---------------------------------------------------------------------\n# 1) \"Umwandlung zu PNAS ...\" paragraph: reorder \"(Andre)\" to the end\n#    and recolor the pieces.\n# ---------------------------------------------------------------------\n$umwandlungPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.IndexOf(\"Umwandlung zu PNAS\") -ge 0) {\n        $umwandlungPara = $p\n        break\n    }\n}\n\nif ($umwandlungPara -ne $null) {\n    $full = $umwandlungPara.Range.Duplicate\n    $full.End = $full.End - 1   # exclude the trailing paragraph mark\n    $full.Text = \"Umwandlung zu PNAS,wenn m\u00f6glich mit Harvard Citation Stil (Andre),\"\n\n    $greenRange = $umwandlungPara.Range.Duplicate\n    $greenRange.Find.ClearFormatting()\n    $greenRange.Find.Text = \"Umwandlung zu PNAS\"\n    $greenRange.Find.Forward = $true\n    $greenRange.Find.Wrap = 1\n    if ($greenRange.Find.Execute()) {\n        $greenRange.Font.Color = $green\n    }\n\n    $amberRange = $umwandlungPara.Range.Duplicate\n    $amberRange.Find.ClearFormatting()\n    $amberRange.Find.Text = \"wenn m\u00f6glich mit Harvard Citation Stil \"\n    $amberRange.Find.Forward = $true\n    $amberRange.Find.Wrap = 1\n    if ($amberRange.Find.Execute()) {\n        $amberRange.Font.Color = $amber\n    }\n}\n\n# ---------------------------------------------------------------------\n# 2) \"Interpretation (Andre)\" paragraph: remove the stray `_GoBack`\n#    bookmark (its text does not change).\n# ---------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# ---------------------------------------------------------------------\n# 3) \"Muss: ...\" paragraph: color two spans green, text unchanged.\n# ---------------------------------------------------------------------\n$mussPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Muss:\")) {\n        $mussPara = $p\n        break\n    }\n}\n\nif ($mussPara -ne $null) {\n    $greenSpan1 = $mussPara.Range.Duplicate\n    $greenSpan1.Find.ClearFormatting()\n    $greenSpan1.Find.Text = \"Titel, Autorenliste, Affiliationen, Abstract (Englisch), Einleitung (\\dropcap) (warum ist es wichtig), \"\n    $greenSpan1.Find.MatchWildcards = $false\n    $greenSpan1.Find.Forward = $true\n    $greenSpan1.Find.Wrap = 1\n    if ($greenSpan1.Find.Execute()) {\n        $greenSpan1.Font.Color = $green\n    }\n\n    $greenSpan2 = $mussPara.Range.Duplicate\n    $greenSpan2.Find.ClearFormatting()\n    $greenSpan2.Find.Text = \"Literaturangabe\"\n    $greenSpan2.Find.Forward = $true\n    $greenSpan2.Find.Wrap = 1\n    if ($greenSpan2.Find.Execute()) {\n        $greenSpan2.Font.Color = $green\n    }\n}\n\n# ---------------------------------------------------------------------\n# 4) Re-insert the `_GoBack` bookmark on the empty paragraph right\n#    after \"Kann: ...\" (Word moves its \"last edit\" bookmark here).\n# ---------------------------------------------------------------------\n$kannIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.StartsWith(\"Kann:\")) {\n        $kannIndex = $i\n        break\n    }\n}\n\nif ($kannIndex -ge 0) {\n    $targetPara = $d.Paragraphs.Item($kannIndex + 1)\n    $d.Bookmarks.Add(\"_GoBack\", $targetPara.Range)\n}\n"}
